$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.310.79"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.912.83"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "349.08"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.38"
$ws.Range("E6").Value = "  -5.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.553"
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.56"
$ws.Range("E10").Value = "  -5.02%  "
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0844"
$ws.Range("E12").Value = "  -3.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.84"
$ws.Range("E13").Value = "  -6.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.372.16"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.56"
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.917.90"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.952"
$ws.Range("E17").Value = "  -3.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.289.37"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.44"
$ws.Range("E19").Value = "  +4.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.34"
$ws.Range("E21").Value = "  -6.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.69"
$ws.Range("E23").Value = "  -3.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "258.94"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.173"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.27"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("E29").Value = "  +5.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.104"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.18"
$ws.Range("E31").Value = "  -4.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.18"
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.42"
$ws.Range("E34").Value = "  -4.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.33"
$ws.Range("E35").Value = "  -5.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0424"
$ws.Range("E37").Value = "  -6.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.11"
$ws.Range("E38").Value = "  -8.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.54"
$ws.Range("E39").Value = "  -6.04%  "
$ws.Range("E40").Value = "  -6.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.30"
$ws.Range("E43").Value = "  -5.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.06"
$ws.Range("E44").Value = "  +7.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.12"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.090.36"
$ws.Range("E46").Value = "  -4.99%  "
$ws.Range("E47").Value = "  -6.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.28"
$ws.Range("E48").Value = "  -9.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.237"
$ws.Range("E49").Value = "  -4.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0333"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.894"
$ws.Range("E51").Value = "  -6.14%  "
